$d = $word.ActiveDocument
$d.Styles.Item("Normal").AutomaticallyUpdate = $true
